$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.919.21'
$ws.Range("D3").Value = '1.632.83'
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.66'
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.70'
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0609'
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.866.59'
$ws.Range("E12").Value = '  +1.92%  '
$ws.Range("D13").Value = '1.633.21'
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.563'
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.28'
$ws.Range("E15").Value = '  +14.85%  '
$ws.Range("D16").Value = '29.943.05'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.60'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").Value = '0.0₃0701'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.84'
$ws.Range("E22").Value = '  +3.61%  '
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.63'
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.48'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0491'
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("E32").Value = '  +4.48%  '
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").Value = '1.428.96'
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.66'
$ws.Range("E35").Value = '  +5.38%  '
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.77'
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.70'
$ws.Range("E40").Value = '  +12.26%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.829'
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0492'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '51.27'
$ws.Range("E47").Value = '  -7.32%  '
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("D49").Value = '1.773.66'
$ws.Range("E49").Value = '  +2.04%  '
$ws.Range("E50").Value = '  +12.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.29'
$ws.Range("E51").Value = '  +4.14%  '
